$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "1.014") must be forced to
# Text format first, otherwise Excel auto-converts them to numbers and
# the original textual formatting (trailing zeros, multi-dot grouping,
# etc.) is lost. NumberFormat is reset back to "Normal" afterwards so no
# stray style index is left on the cell.
$textForceCells = @(
    "D4",
    "D5",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D17",
    "D18",
    "D19",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.557.73"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.857.61"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +0.98%  "
$ws.Range("D5").Value = "332.93"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").Value = "0.4659"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("D8").Value = "0.3879"
$ws.Range("E8").Value = "  -2.07%  "
$ws.Range("D9").Value = "45.70"
$ws.Range("E9").Value = "  -4.88%  "
$ws.Range("D10").Value = "0.07956"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "0.9933"
$ws.Range("E11").Value = "  -3.94%  "
$ws.Range("D12").Value = "21.48"
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").Value = "1.860.85"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "5.955"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "7.176"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "87.55"
$ws.Range("D18").Value = "0.06709"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "0.00001039"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "27.554.27"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "5.433"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").Value = "10.79"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").Value = "2.315"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").Value = "2.081.31"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "158.28"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "19.65"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "2.095"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "5.342"
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "121.10"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "0.9648"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("D33").Value = "0.09445"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").Value = "3.641"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "5.267"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").Value = "1.313"
$ws.Range("E36").Value = "  -9.58%  "
$ws.Range("D37").Value = "0.06004"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("D38").Value = "0.02208"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("D39").Value = "1.192"
$ws.Range("E39").Value = "  -3.72%  "
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "8.106"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").Value = "0.5876"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").Value = "0.1869"
$ws.Range("D44").Value = "10.14"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("D45").Value = "1.252"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").Value = "0.5591"
$ws.Range("E46").Value = "  -2.70%  "
$ws.Range("D47").Value = "12.00"
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("D48").Value = "1.903"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").Value = "3.272"
$ws.Range("E49").Value = "  -3.06%  "
$ws.Range("D50").Value = "0.06749"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").Value = "111.93"
$ws.Range("E51").Value = "  -1.94%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
